# Add diary entries for week 7 (rows 26-29) on the existing worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 26-29 are currently blank template rows (style-only). Row 23 already
# carries the exact target style combination for columns A-G
# (A=18, B=22, C..G=20), so copy its formatting down into the four new rows
# before writing values, matching the existing "filled-in" weeks above them.
$ws.Range("A23:G23").Copy() | Out-Null
$ws.Range("A26:G26").PasteSpecial(-4122) | Out-Null
$ws.Range("A27:G27").PasteSpecial(-4122) | Out-Null
$ws.Range("A28:G28").PasteSpecial(-4122) | Out-Null
$ws.Range("A29:G29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 26 - Tue 2020-02-20 (serial 43881)
$ws.Range("A26").Value = 43881
$ws.Range("B26").Value = "5:00PM - 8:00PM"
$ws.Range("C26").Value = "N/A"
$ws.Range("D26").Value = "Learn something from class"
$ws.Range("E26").Value = "Learned about how to abstract out the architecture from a system"
$ws.Range("F26").Value = "Runelite should be fairly easy to do for us as its structure is very well defined."
$ws.Range("G26").Value = "Decent"

# Row 27 - Sat 2020-02-24 (serial 43885)
$ws.Range("A27").Value = 43885
$ws.Range("B27").Value = "4:00PM - 6:00PM"
$ws.Range("C27").Value = "Thuc, Harry"
$ws.Range("D27").Value = "Find and document architeture for runelite"
$ws.Range("E27").Value = "Was able to contact the devs directly and get advice on generating architecture diagram"
$ws.Range("F27").Value = "Was surprised the devs are this responsive on their official discord server, especially the creator of the project being able to directly answer my questions felt great."
$ws.Range("G27").Value = "Great"

# Row 28 - Sun 2020-02-25 (serial 43886)
$ws.Range("A28").Value = 43886
$ws.Range("B28").Value = "9:00PM - 1:00AM"
$ws.Range("C28").Value = "Thuc, Harry"
$ws.Range("D28").Value = "Find and document social context and contribution guidelines for runelite"
$ws.Range("E28").Value = "Was able to scrape a lot of info from their github page, and official website."
$ws.Range("F28").Value = "Having a well maintained project that is still ongoing helps alot with getting such info, especially when I have been using it for over 3 years"
$ws.Range("G28").Value = "Great"

# Row 29 - Mon 2020-02-26 (serial 43887)
$ws.Range("A29").Value = 43887
$ws.Range("B29").Value = "9:00PM - 5:00AM"
$ws.Range("C29").Value = "Thuc, Harry"
$ws.Range("D29").Value = "Find interesting pull requests and issues, and document them. Reformat the document."
$ws.Range("E29").Value = "Found a lot of interesting stuff happening in the official github."
$ws.Range("F29").Value = "Some issues can only be solved by the community, some pull requests were denied and closed, learned a lot of what the maintainers are expecting from each contribution"
$ws.Range("G29").Value = "exhausted"
